$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (shifts old rows 4-30 down to 5-31),
# inheriting the formatting of the row above (row 3) like Excel's default Insert behavior.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the "SpecimenAntibodyResults7" test record.
$ws.Range("A4").Value = "SpecimenAntibodyResults7"
$ws.Range("B4").Value = "G814450907"
$ws.Range("C4").Value = "Non-Negative"
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = "C07"
$ws.Range("F4").Value = "20200617_1"
$ws.Range("G4").Value = "AntibodyResults1"
$ws.Range("H4").Value = "Rack 02"
# Column I keeps its quote-prefixed text style across the whole sheet; a leading
# apostrophe preserves that formatting (matching the rest of the column) instead
# of resetting it the way a plain assignment would.
$ws.Range("I4").Value = "'RackPos 04"

# Move the active selection to E5, matching the saved workbook view.
$ws.Range("E5").Select() | Out-Null
